$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 637941.8
$ws.Range("I28").Value = 865313.9
$ws.Range("K28").Value = 865313.9
$ws.Range("M28").Value = -864828.9

$ws.Range("H47").Value = 26000
$ws.Range("I47").Value = 20000
$ws.Range("J47").Value = 29000
$ws.Range("K47").Value = 20000
$ws.Range("L47").Value = 29000
$ws.Range("M47").Value = -19028
$ws.Range("N47").Value = -30944

$ws.Range("H48").Value = 6055.5
$ws.Range("I48").Value = 6055.5
$ws.Range("K48").Value = 18166.5
$ws.Range("M48").Value = -17874.5

$ws.Range("H56").Value = 6055.5
$ws.Range("I56").Value = 6055.5
$ws.Range("K56").Value = 18166.5
$ws.Range("M56").Value = -17632.5

$ws.Range("H106").Value = 7008608
$ws.Range("I106").Value = 8625426
$ws.Range("K106").Value = 8625426
$ws.Range("M106").Value = -8624795

$ws.Range("H137").Value = 43479464
$ws.Range("I137").Value = 76924070
$ws.Range("J137").Value = 1470
$ws.Range("K137").Value = 230772210
$ws.Range("L137").Value = 4410
$ws.Range("M137").Value = -230769660
$ws.Range("N137").Value = -9510

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4122.5
$ws.Range("I2").Value = 5675
$ws.Range("J2").Value = 3346.25
$ws.Range("K2").Value = 5675
$ws.Range("L2").Value = 3346.25
$ws.Range("M2").Value = -5562
$ws.Range("N2").Value = -3572.25

$ws.Range("H32").Value = 22109.037
$ws.Range("I32").Value = 2095.8157
$ws.Range("J32").Value = 66844.47
$ws.Range("K32").Value = 2095.8157
$ws.Range("L32").Value = 66844.47
$ws.Range("M32").Value = -1808.8157
$ws.Range("N32").Value = -67418.47

$ws.Range("H61").Value = 2523.2285
$ws.Range("I61").Value = 1502.9584
$ws.Range("J61").Value = 4749.273
$ws.Range("K61").Value = 1502.9584
$ws.Range("L61").Value = 4749.273
$ws.Range("M61").Value = -1290.9584
$ws.Range("N61").Value = -5173.273

$ws.Range("H116").Value = 4122.5
$ws.Range("I116").Value = 5675
$ws.Range("J116").Value = 3346.25
$ws.Range("K116").Value = 5675
$ws.Range("L116").Value = 3346.25
$ws.Range("M116").Value = -3381
$ws.Range("N116").Value = -7934.25

$ws.Range("H122").Value = 1817.0667
$ws.Range("I122").Value = 1817.0667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5451.2001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3001.2001
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 1833.15
$ws.Range("I132").Value = 1485.0834
$ws.Range("J132").Value = 4965.75
$ws.Range("K132").Value = 4455.2502
$ws.Range("L132").Value = 14897.25
$ws.Range("M132").Value = -1925.2502
$ws.Range("N132").Value = -19957.25

$ws.Range("H136").Value = 2523.2285
$ws.Range("I136").Value = 1502.9584
$ws.Range("J136").Value = 4749.273
$ws.Range("K136").Value = 4508.8752
$ws.Range("L136").Value = 14247.819
$ws.Range("M136").Value = -1958.8752
$ws.Range("N136").Value = -19347.819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4122.5
$ws.Range("I3").Value = 5675
$ws.Range("J3").Value = 3346.25
$ws.Range("K3").Value = 5675
$ws.Range("L3").Value = 3346.25
$ws.Range("M3").Value = -5561
$ws.Range("N3").Value = -3574.25

$ws.Range("H105").Value = 9105.111000000001
$ws.Range("I105").Value = 11084.154
$ws.Range("K105").Value = 11084.154
$ws.Range("M105").Value = -9337.154

$ws.Range("H134").Value = 3208.0303
$ws.Range("I134").Value = 2310.6191
$ws.Range("J134").Value = 4778.5
$ws.Range("K134").Value = 6931.8573
$ws.Range("L134").Value = 14335.5
$ws.Range("M134").Value = -4396.8573
$ws.Range("N134").Value = -19405.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1499.421
$ws.Range("I31").Value = 1417.5
$ws.Range("J31").Value = 1936.3334
$ws.Range("K31").Value = 1417.5
$ws.Range("L31").Value = 1936.3334
$ws.Range("M31").Value = -1122.5
$ws.Range("N31").Value = -2526.3334

$ws.Range("H34").Value = 1499.421
$ws.Range("I34").Value = 1417.5
$ws.Range("J34").Value = 1936.3334
$ws.Range("K34").Value = 1417.5
$ws.Range("L34").Value = 1936.3334
$ws.Range("M34").Value = -1215.5
$ws.Range("N34").Value = -2340.3334

$ws.Range("H59").Value = 37500
$ws.Range("J59").Value = 37500
$ws.Range("L59").Value = 37500
$ws.Range("N59").Value = -39790

$ws.Range("H107").Value = 743.36365
$ws.Range("I107").Value = 440
$ws.Range("K107").Value = 440
$ws.Range("M107").Value = 1480

$ws.Range("H132").Value = 2830.9697
$ws.Range("I132").Value = 2026.909
$ws.Range("J132").Value = 4439.091
$ws.Range("K132").Value = 6080.727000000001
$ws.Range("L132").Value = 13317.273
$ws.Range("M132").Value = -3550.727000000001
$ws.Range("N132").Value = -18377.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4761.5386
$ws.Range("I56").Value = 4761.5386
$ws.Range("K56").Value = 4761.5386
$ws.Range("M56").Value = -4231.5386

$ws.Range("H131").Value = 1618.44
$ws.Range("J131").Value = 1894.975
$ws.Range("L131").Value = 5684.924999999999
$ws.Range("N131").Value = -15764.925

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1173.1333
$ws.Range("I102").Value = 970.8182
$ws.Range("K102").Value = 970.8182
$ws.Range("M102").Value = 651.1818

$ws.Range("H126").Value = 3109.5
$ws.Range("I126").Value = 2480.5
$ws.Range("J126").Value = 3361.1
$ws.Range("K126").Value = 7441.5
$ws.Range("L126").Value = 10083.3
$ws.Range("M126").Value = -4971.5
$ws.Range("N126").Value = -15023.3

$ws.Range("H132").Value = 2066.1052
$ws.Range("I132").Value = 1506.32
$ws.Range("J132").Value = 3142.6155
$ws.Range("K132").Value = 4518.96
$ws.Range("L132").Value = 9427.8465
$ws.Range("M132").Value = -1988.96
$ws.Range("N132").Value = -14487.8465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3470.5881
$ws.Range("I40").Value = 3250
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 3250
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -3114
$ws.Range("N40").Value = -3772

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H122").Value = 3334.682
$ws.Range("I122").Value = 2163.1428
$ws.Range("J122").Value = 3881.4
$ws.Range("K122").Value = 6489.428400000001
$ws.Range("L122").Value = 11644.2
$ws.Range("M122").Value = -4039.428400000001
$ws.Range("N122").Value = -16544.2

$ws.Range("H136").Value = 5335.385
$ws.Range("I136").Value = 3067.5386
$ws.Range("J136").Value = 7603.231
$ws.Range("K136").Value = 9202.6158
$ws.Range("L136").Value = 22809.693
$ws.Range("M136").Value = -6652.6158
$ws.Range("N136").Value = -27909.693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 17297.777
$ws.Range("I62").Value = 26158.4
$ws.Range("J62").Value = 6222
$ws.Range("K62").Value = 26158.4
$ws.Range("L62").Value = 6222
$ws.Range("M62").Value = -25534.4
$ws.Range("N62").Value = -7470

$ws.Range("H65").Value = 17297.777
$ws.Range("I65").Value = 26158.4
$ws.Range("J65").Value = 6222
$ws.Range("K65").Value = 130792
$ws.Range("L65").Value = 31110
$ws.Range("M65").Value = -127672
$ws.Range("N65").Value = -37350

$ws.Range("H113").Value = 325
$ws.Range("I113").Value = 233.33333
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 699.99999
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 1470.00001
$ws.Range("N113").Value = -6140

$ws.Range("H132").Value = 23813132
$ws.Range("I132").Value = 38464524
$ws.Range("J132").Value = 4621.875
$ws.Range("K132").Value = 115393572
$ws.Range("L132").Value = 13865.625
$ws.Range("M132").Value = -115391042
$ws.Range("N132").Value = -18925.625
